# "Add bottom of screen back/next"
#
# - settings!B3 (form_version) bumped to the new build id
# - settings gains a new "showFooter" = 1 row so the generated form shows
#   the back/next controls at the bottom of the screen
# - table_specific_translations: a couple of cells that had picked up a
#   slightly-different (but visually identical) Arial font are normalised
#   back to the sheet's standard style
# - settings becomes the active sheet/tab, with B4 selected

$wb = $excel.ActiveWorkbook

# --- settings sheet -------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

# Bump form_version
$settings.Range("B3").Value = 20210221001

# Add a new setting row: showFooter = 1
$settings.Cells.Item(10, 1).Value = "showFooter"
$settings.Cells.Item(10, 2).Value = 1

# --- table_specific_translations sheet -------------------------------
$translations = $wb.Worksheets.Item("table_specific_translations")

# Normalise B4/C4 and C11/D11 to the same font as the rest of the column
# (drops the stray duplicate font/style that only differed by a missing
# charset attribute).
$translations.Range("B4").Font.Name = "Arial"
$translations.Range("B4").Font.Size = 10
$translations.Range("C4").Font.Name = "Arial"
$translations.Range("C4").Font.Size = 10
$translations.Range("C11").Font.Name = "Arial"
$translations.Range("C11").Font.Size = 10
$translations.Range("D11").Font.Name = "Arial"
$translations.Range("D11").Font.Size = 10

# --- active sheet / selection ----------------------------------------
$settings.Activate()
$settings.Range("B4").Select()
